# Auto update Excel log
# Appends newly-logged sensor rows to the PIR, Humidity and Temperature sheets,
# mirroring the continuous sensor log append pattern already used throughout
# this workbook.

$wb = $excel.ActiveWorkbook

function Append-Rows {
    param(
        $SheetName,
        $StartRow,
        $Rows,
        $TextCols
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $endRow = $StartRow + $Rows.Count - 1

    # Some columns hold values that look numeric/date-like to Excel's
    # automatic type detection (e.g. "2026-01-28" looks like a date,
    # "87.7%" looks like a percentage). Assigning such text via .Value
    # would silently convert it to a date-serial / fractional number with
    # an applied number format. Force those destination columns to Text
    # first (using string concatenation, not interpolation, for the
    # address - interpolation of function parameters into a Range address
    # has proven unreliable in this runtime)...
    $textRanges = @()
    foreach ($col in $TextCols) {
        $addr = $col + $StartRow + ":" + $col + $endRow
        $rng = $ws.Range($addr)
        $rng.NumberFormat = "@"
        $textRanges += $rng
    }

    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $StartRow + $i
        $row = $Rows[$i]
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]
    }

    # ...then strip the temporary formatting back off again so the new
    # cells stay styled exactly like every other cell in the log (no
    # explicit cell style / number format applied).
    foreach ($rng in $textRanges) {
        $rng.ClearFormats()
    }
}

# ---- PIR sheet: 5 new rows (172-176) ----
$pirRows = @(
    ,@("2026-01-28", "17:29:41", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "17:29:42", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "17:29:45", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "17:29:51", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "17:29:55", "17:00", "Bathroom", "No Motion", "Inactive")
)
Append-Rows "PIR" 172 $pirRows @("A")

# ---- Humidity sheet: 6 new rows (172-177) ----
$humidityRows = @(
    ,@("2026-01-28", "17:29:38", "17:00", "Bathroom", "87.7%", "Active")
    ,@("2026-01-28", "17:29:40", "17:00", "Bathroom", "87.7%", "Active")
    ,@("2026-01-28", "17:29:41", "17:00", "Bathroom", "87.7%", "Active")
    ,@("2026-01-28", "17:29:43", "17:00", "Bathroom", "87.7%", "Active")
    ,@("2026-01-28", "17:29:50", "17:00", "Bathroom", "87.7%", "Active")
    ,@("2026-01-28", "17:29:58", "17:00", "Bathroom", "87.7%", "Active")
)
Append-Rows "Humidity" 172 $humidityRows @("A", "E")

# ---- Temperature sheet: 6 new rows (172-177) ----
$temperatureRows = @(
    ,@("2026-01-28", "17:29:38", "17:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "17:29:40", "17:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "17:29:42", "17:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "17:29:44", "17:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "17:29:50", "17:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "17:29:58", "17:00", "Bathroom", "22.8C", "Active")
)
Append-Rows "Temperature" 172 $temperatureRows @("A")
